# Updated cryptos list on Mon Feb 19 14:54:00 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price ("D") column cells below look like plain decimal numbers, so Excel's
# Range.Value setter would silently coerce them to numeric cells (dropping
# e.g. a trailing ".70" -> "7.7"). Force these specific cells to keep a Text
# number format first so the literal string is preserved, matching the
# original sheet where every Price cell is stored as text.
$textPriceCells = @(
    "D5","D6","D10","D11","D12","D14","D17","D21","D23","D24","D25","D26",
    "D27","D31","D32","D33","D34","D35","D38","D39","D41","D43","D44","D46",
    "D49","D50"
)
foreach ($addr in $textPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "52.153.16"
$ws.Range("E2").Value = "  +0.49%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.904.81"
$ws.Range("E3").Value = "  +3.42%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.06%  "

# Row 5 - BNB
$ws.Range("D5").Value = "352.01"
$ws.Range("E5").Value = "  -0.60%  "

# Row 6 - Solana
$ws.Range("D6").Value = "112.91"
$ws.Range("E6").Value = "  +0.62%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  -0.37%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.03%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -1.21%  "

# Row 10 - Avalanche
$ws.Range("D10").Value = "39.57"
$ws.Range("E10").Value = "  -2.04%  "

# Row 11 - was TRON, becomes Dogecoin
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").Value = "0.0863"
$ws.Range("E11").Value = "  +2.83%  "

# Row 12 - was Dogecoin, becomes TRON
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "0.136"
$ws.Range("E12").Value = "  +0.58%  "

# Row 13 - Chainlink
$ws.Range("E13").Value = "  -1.22%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "7.70"
$ws.Range("E14").Value = "  -1.18%  "

# Row 15 - Wrapped liquid staked Ether 2.0
$ws.Range("D15").Value = "3.359.52"
$ws.Range("E15").Value = "  +3.47%  "

# Row 16 - Wrapped Ether
$ws.Range("D16").Value = "2.918.26"
$ws.Range("E16").Value = "  +3.90%  "

# Row 17 - Polygon
$ws.Range("D17").Value = "0.983"
$ws.Range("E17").Value = "  +3.81%  "

# Row 18 - Wrapped BTC
$ws.Range("D18").Value = "52.188.04"
$ws.Range("E18").Value = "  +0.64%  "

# Row 19 - ImmutableX
$ws.Range("E19").Value = "  +2.56%  "

# Row 20 - Uniswap
$ws.Range("E20").Value = "  -0.59%  "

# Row 21 - Internet Computer (DFINITY)
$ws.Range("D21").Value = "13.85"
$ws.Range("E21").Value = "  +1.66%  "

# Row 22 - Shiba Inu
$ws.Range("D22").Value = "0.0₃0973"
$ws.Range("E22").Value = "  -0.14%  "

# Row 23 - Litecoin
$ws.Range("D23").Value = "71.11"
$ws.Range("E23").Value = "  +1.11%  "

# Row 24 - Bitcoin Cash
$ws.Range("D24").Value = "269.25"
$ws.Range("E24").Value = "  +0.50%  "

# Row 25 - PancakeSwap
$ws.Range("D25").Value = "2.79"
$ws.Range("E25").Value = "  +1.12%  "

# Row 26 - Kaspa
$ws.Range("D26").Value = "0.181"
$ws.Range("E26").Value = "  +13.46%  "

# Row 27 - Ethereum Classic
$ws.Range("D27").Value = "26.74"
$ws.Range("E27").Value = "  +2.06%  "

# Row 28 - Dai
$ws.Range("E28").Value = "  -0.09%  "

# Row 29 - Cosmos
$ws.Range("E29").Value = "  +2.29%  "

# Row 30 - Hedera
$ws.Range("E30").Value = "  +15.69%  "

# Row 31 - Filecoin
$ws.Range("D31").Value = "6.64"
$ws.Range("E31").Value = "  +8.08%  "

# Row 32 - Injective Protocol
$ws.Range("D32").Value = "37.39"
$ws.Range("E32").Value = "  -4.17%  "

# Row 33 - Toncoin
$ws.Range("D33").Value = "2.26"
$ws.Range("E33").Value = "  -0.74%  "

# Row 34 - Render Token
$ws.Range("D34").Value = "6.19"
$ws.Range("E34").Value = "  +11.62%  "

# Row 35 - OKB
$ws.Range("D35").Value = "53.05"
$ws.Range("E35").Value = "  +0.97%  "

# Row 36 - VeChain
$ws.Range("E36").Value = "  -1.51%  "

# Row 37 - First Digital USD
$ws.Range("E37").Value = "  -0.16%  "

# Row 38 - Lido DAO Token
$ws.Range("D38").Value = "3.30"
$ws.Range("E38").Value = "  +3.98%  "

# Row 39 - Celestia
$ws.Range("D39").Value = "18.79"
$ws.Range("E39").Value = "  -0.78%  "

# Row 40 - ARBITRUM
$ws.Range("E40").Value = "  +1.23%  "

# Row 41 - Stacks
$ws.Range("D41").Value = "2.72"
$ws.Range("E41").Value = "  +7.31%  "

# Row 42 - Stellar
$ws.Range("E42").Value = "  +1.24%  "

# Row 43 - EnergySwap
$ws.Range("D43").Value = "22.95"
$ws.Range("E43").Value = "  +4.74%  "

# Row 44 - Monero
$ws.Range("D44").Value = "119.24"
$ws.Range("E44").Value = "  -0.80%  "

# Row 45 - WEMIX Token
$ws.Range("E45").Value = "  -1.68%  "

# Row 46 - ApeX Protocol
$ws.Range("D46").Value = "2.56"
$ws.Range("E46").Value = "  +3.98%  "

# Row 47 - Maker
$ws.Range("D47").Value = "2.170.40"
$ws.Range("E47").Value = "  +3.05%  "

# Row 48 - NEAR Protocol
$ws.Range("E48").Value = "  -0.51%  "

# Row 49 - The Graph
$ws.Range("D49").Value = "0.264"
$ws.Range("E49").Value = "  +20.80%  "

# Row 50 - BEAM
$ws.Range("D50").Value = "0.0341"
$ws.Range("E50").Value = "  +9.26%  "

# Row 51 - SEI
$ws.Range("E51").Value = "  -1.06%  "
